$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (shifts existing rows 12-34 down to 13-35)
$ws.Rows.Item(12).Insert()

# Populate the new "Transport" row (row 12) - simple blade transport model
$ws.Range("A12").Value2 = "WISDEM"
$ws.Range("B12").Value2 = 100
$ws.Range("C12").Value2 = 1.5
$ws.Range("D12").Value2 = "TransportCost"
$ws.Range("E12").Value2 = "Transport"
$ws.Range("F12").Value2 = "Other"
$ws.Range("G12").Value2 = 1190000
$ws.Range("H12").Style = "Normal"

# Crane edits: adjust a few downstream Management cost figures
$ws.Range("G29").Value2 = 197951.126780675
$ws.Range("G32").Value2 = 353484.15496549098
$ws.Range("G33").Value2 = 4602363.6976506999

# Refresh the AutoFilter range to cover the new row count
$ws.AutoFilterMode = $false
$ws.Range("A1:G35").AutoFilter()

# Update the _FilterDatabase defined name to match the new filter range
$n = $wb.Names.Item(1)
$n.RefersTo = "=costs_by_module_type_operation!`$A`$1:`$G`$35"

# Update the view: select I12, reset the frozen/top-left scroll position
$ws.Activate()
$ws.Range("I12").Select()
